$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 3 through 11 (old data rows beyond the new single data row)
$ws.Range("A3:E11").Delete()

# Update header row
$ws.Range("A1").Value = "Qtd_Nós"
$ws.Range("B1").Value = "Ativos"
$ws.Range("C1").Value = "Distancia"
$ws.Range("D1").Value = "Tempo"

# Clear old column E header / leftover data
$ws.Range("E1").ClearContents()

# Update data row 2
$ws.Range("A2").Value = 81
$ws.Range("B2").Value = 30
$ws.Range("C2").Value = 8837
$ws.Range("D2").Value = 0.1609594821929932
$ws.Range("E2").ClearContents()
